$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 60.87

# Clear the Scope ID # value (was "#NO MATCH", now blank)
$ws.Range("G10").Value = ""

# Update line item pricing and total
$ws.Range("H16").Value = 60.87
$ws.Range("H17").Value = 60.87
